# Update the two-digit multiplication problems in the single table.
# Each data row of the table (rows 1, 5, 10, 15, 20) holds 5 cells with
# expressions like "NN×NN=". We set each cell's Range.Text directly
# (rather than using document-wide Find/Replace) so that duplicate
# expressions (e.g. "72×58=" appears in two different cells) are updated
# independently with their correct new values, and existing run
# formatting on the cell is preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; New="27×79="},
    @{Row=1;  Col=2; New="56×38="},
    @{Row=1;  Col=3; New="28×84="},
    @{Row=1;  Col=4; New="65×59="},
    @{Row=1;  Col=5; New="87×32="},

    @{Row=5;  Col=1; New="11×98="},
    @{Row=5;  Col=2; New="99×26="},
    @{Row=5;  Col=3; New="24×49="},
    @{Row=5;  Col=4; New="28×68="},
    @{Row=5;  Col=5; New="14×91="},

    @{Row=10; Col=1; New="29×64="},
    @{Row=10; Col=2; New="67×53="},
    @{Row=10; Col=3; New="30×16="},
    @{Row=10; Col=4; New="55×11="},
    @{Row=10; Col=5; New="68×83="},

    @{Row=15; Col=1; New="24×11="},
    @{Row=15; Col=2; New="98×64="},
    @{Row=15; Col=3; New="11×89="},
    @{Row=15; Col=4; New="52×16="},
    @{Row=15; Col=5; New="56×23="},

    @{Row=20; Col=1; New="72×74="},
    @{Row=20; Col=2; New="70×34="},
    @{Row=20; Col=3; New="19×57="},
    @{Row=20; Col=4; New="28×37="},
    @{Row=20; Col=5; New="77×47="}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cell.Range.Text = $item.New
}

Write-Host "Replacements applied:" $replacements.Count
